$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.498.52'
$ws.Range('E2').Value = '  -0.58%  '

# Row 3
$ws.Range('D3').Value = '1.617.40'
$ws.Range('E3').Value = '  -1.57%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.21'
$ws.Range('E5').Value = '  -0.87%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.523'
$ws.Range('E6').Value = '  -1.29%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.74'
$ws.Range('E8').Value = '  -1.10%  '

# Row 9
$ws.Range('E9').Value = '  +2.10%  '

# Row 10
$ws.Range('E10').Value = '  +0.16%  '

# Row 11
$ws.Range('E11').Value = '  -0.46%  '

# Row 12
$ws.Range('D12').Value = '1.847.09'
$ws.Range('E12').Value = '  -1.69%  '

# Row 13
$ws.Range('D13').Value = '1.616.14'
$ws.Range('E13').Value = '  -1.60%  '

# Row 14
$ws.Range('E14').Value = '  -0.50%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.549'
$ws.Range('E15').Value = '  -2.17%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.96'
$ws.Range('E16').Value = '  +1.29%  '

# Row 17
$ws.Range('D17').Value = '27.491.39'
$ws.Range('E17').Value = '  -0.34%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.70'
$ws.Range('E18').Value = '  +0.24%  '

# Row 19
$ws.Range('E19').Value = '  -0.67%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.52'
$ws.Range('E20').Value = '  -1.45%  '

# Row 21
$ws.Range('E21').Value = '  +0.01%  '

# Row 22
$ws.Range('E22').Value = '  -0.69%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.17'
$ws.Range('E23').Value = '  +1.15%  '

# Row 24
$ws.Range('E24').Value = '  +6.52%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.67'
$ws.Range('E25').Value = '  +0.33%  '

# Row 26
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.111'
$ws.Range('E26').Value = '  -1.31%  '

# Row 27
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.03%  '

# Row 28
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.81'
$ws.Range('E28').Value = '  -1.99%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.54'
$ws.Range('E29').Value = '  -0.59%  '

# Row 30
$ws.Range('E30').Value = '  -0.55%  '

# Row 31
$ws.Range('E31').Value = '  -0.50%  '

# Row 32
$ws.Range('E32').Value = '  -0.99%  '

# Row 33
$ws.Range('D33').Value = '1.441.05'
$ws.Range('E33').Value = '  +0.22%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.06'
$ws.Range('E34').Value = '  -3.26%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').Value = '  -3.54%  '

# Row 36
$ws.Range('E36').Value = '  -0.49%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.938'
$ws.Range('E37').Value = '  +4.52%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.562'
$ws.Range('E38').Value = '  -1.61%  '

# Row 39
$ws.Range('E39').Value = '  -0.16%  '

# Row 40
$ws.Range('E40').Value = '  -2.22%  '

# Row 41
$ws.Range('E41').Value = '  -0.09%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '67.74'
$ws.Range('E42').Value = '  +3.74%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.991'
$ws.Range('E43').Value = '  -3.97%  '

# Row 44
$ws.Range('E44').Value = '  -0.73%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.40'
$ws.Range('E45').Value = '  -5.21%  '

# Row 46
$ws.Range('E46').Value = '  -2.14%  '

# Row 47
$ws.Range('D47').Value = '1.757.90'
$ws.Range('E47').Value = '  -1.75%  '

# Row 48
$ws.Range('E48').Value = '  +0.81%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.39'
$ws.Range('E49').Value = '  -0.02%  '

# Row 50
$ws.Range('E50').Value = '  +2.65%  '

# Row 51
$ws.Range('E51').Value = '  +1.64%  '
